# "Add some PT Server"
# Adds three new rows to the server list for Portugal (PT):
#   - DataPacket (Lisboa)
#   - NOS Lisboa
#   - NOS Porto
#
# The rows are inserted in the same order the original author appears to
# have entered them (NOS Lisboa + NOS Porto first, then DataPacket above
# them) so that the shared-string table ends up in the same order as the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two rows right after the last existing row (122) for
#     "NOS Lisboa" and "NOS Porto", copying the formatting (borders, etc.)
#     from the row directly above so the new cells keep style index 2.
$ws.Rows("123:124").Insert()
$ws.Range("A122:E122").Copy()
$ws.Range("A123:E124").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A123").Value = "PT"
$ws.Range("B123").Value = "NOS Lisboa"
$ws.Range("C123").Value = "10Gbps"
$ws.Range("D123").Value = "NOS_10Gbps_Lisboa_PT.xml"
$ws.Range("E123").Value = "No"

$ws.Range("A124").Value = "PT"
$ws.Range("B124").Value = "NOS Porto"
$ws.Range("C124").Value = "10Gbps"
$ws.Range("D124").Value = "NOS_10Gbps_Porto_PT.xml"
$ws.Range("E124").Value = "No"

# --- Step 2: insert one more row above those two for "DataPacket".
$ws.Rows("123:123").Insert()
$ws.Range("A124:E124").Copy()
$ws.Range("A123:E123").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A123").Value = "PT"
$ws.Range("B123").Value = "DataPacket"
$ws.Range("C123").Value = "10Gbps"
$ws.Range("D123").Value = "DataPacket_10Gbps_Lisboa_PT.xml"
$ws.Range("E123").Value = "No"

# --- Leave the selection on one of the newly added rows, similar to where
#     the author's cursor ended up after the edit.
$ws.Range("E124").Select()

$excel.CutCopyMode = 0
